$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.438.46'
$ws.Range('D3').Value = '1.574.55'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '291.31'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +2.14%  '
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3422'
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07674'
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.003'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.30'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.005'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.930'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '1.574.42'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.46'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06775'
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.82'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.237'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.05'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('B24').Value = 'WrappedBTC'
$ws.Range('C24').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D24').Value = '22.430.63'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.423'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.749'
$ws.Range('E26').Value = '  -6.80%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.33'
$ws.Range('E27').Value = '  +2.51%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '145.63'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.045'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '126.15'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('B31').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C31').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D31').Value = '1.753.70'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.230'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.014'
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.013'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '10.04'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08578'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02565'
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2318'
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.347'
$ws.Range('E39').Value = '  +8.27%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06589'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.470'
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6477'
$ws.Range('E42').Value = '  +1.39%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.60'
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.13'
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.002'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6040'
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.791'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.300'
$ws.Range('E48').Value = '  +9.31%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.097'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '125.79'
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07341'
$ws.Range('E51').Value = '  +1.05%  '
